# Weekly price-sheet update: a new record is inserted as row 356 (pushing
# every existing record for rows 356-415 down by one, to 357-416), with the
# sheet's used range growing from A1:R415 to A1:R416.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 356; everything below shifts down.
$ws.Rows.Item(356).Insert()

# Populate the newly inserted row with the new market observation.
$ws.Range("A356").Value = 4
$ws.Range("B356").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C356").Value = "Los Lagos"
$ws.Range("D356").Value = 45218
$ws.Range("E356").Value = 10
$ws.Range("F356").Value = 100112039
$ws.Range("G356").Value = "Ciboulette"
$ws.Range("H356").Value = "Sin especificar"
$ws.Range("I356").Value = "Primera"
$ws.Range("J356").Value = 80
$ws.Range("K356").Value = 3500
$ws.Range("L356").Value = 3500
$ws.Range("M356").Value = 3500
$ws.Range("N356").Value = "$/docena de atados"
$ws.Range("O356").Value = "Región Metropolitana"
$ws.Range("P356").Value = 1167
$ws.Range("Q356").Value = 3
$ws.Range("R356").Value = "Hortaliza"
